# Insert a new data row for "Macroferia Regional de Talca - Ajo" above the
# current row 267 (shifting existing rows 267-356 down to 268-357), then
# populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267, shifting rows 267:356 down to 268:357.
$ws.Rows("267:267").Insert()

# Populate the new row 267 with the new record's data.
$ws.Range("A267").Value = 5
$ws.Range("B267").Value = "Macroferia Regional de Talca"
$ws.Range("C267").Value = "Maule"
$ws.Range("D267").Value = 44809
$ws.Range("D267").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E267").Value = 7
$ws.Range("F267").Value = 100112003
$ws.Range("G267").Value = "Ajo"
$ws.Range("H267").Value = "Chino"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 250
$ws.Range("K267").Value = 25000
$ws.Range("L267").Value = 25000
$ws.Range("M267").Value = 25000
$ws.Range("N267").Value = "`$/malla 10 kilos"
$ws.Range("O267").Value = "China"
$ws.Range("P267").Value = 2500
$ws.Range("Q267").Value = 10
$ws.Range("R267").Value = "Hortaliza"
